# The workbook originally has a single worksheet named "Sayfa1" that holds
# a multiplication table (1..10 x 1..10). The commit renames it to
# "Practice" (part of a larger change adding extra practice/solution
# sheets elsewhere in the project). No cell values/layout change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Practice"
